$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A4:D4")
$rng.NumberFormat = "@"

$ws.Range("A4").Value = "23"
$ws.Range("B4").Value = "2024-09-12"
$ws.Range("C4").Value = "96"
$ws.Range("D4").Value = "12"
